# "Generate Report for Archive"
#
# The localization-status report is regenerated: rows that were previously
# "Ready for handoff" are now "In Translation", and the (zh-cn / de-de /
# Status) columns that carry that text are narrowed to match the shorter
# label.

$wb = $excel.ActiveWorkbook

# --- 1. Update status text wherever it appears -----------------------------
# "Overview" sheet: zh-cn (E) and de-de (F) status columns, rows 2-3.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

# "zh-cn" / "de-de" detail sheets: Status column (C), rows 2-3.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- 2. Narrow the columns that held the old, longer status text -----------
# Same target display width on all four columns (Overview!E:F, zh-cn!C,
# de-de!C) - set via the standard character-based ColumnWidth property.
# (Columns must be addressed by number here - this host's Columns.Item
# does not resolve letter keys like "E".)
$newWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth

Write-Output "Updated status text and column widths for archive report."
